{"js": "// Replace the division problems in the table with their new values.\n// Each old value is unique in the document, so a simple search &\n// replace keyed on the full \"NNN\u00f7N=\" text is unambiguous.\nconst replacements = [\n  [\"919\u00f76=\", \"512\u00f75=\"],\n  [\"585\u00f79=\", \"397\u00f76=\"],\n  [\"582\u00f72=\", \"685\u00f74=\"],\n  [\"372\u00f74=\", \"141\u00f77=\"],\n  [\"843\u00f74=\", \"697\u00f79=\"],\n  [\"983\u00f74=\", \"193\u00f79=\"],\n  [\"964\u00f76=\", \"536\u00f72=\"],\n  [\"583\u00f73=\", \"653\u00f78=\"],\n  [\"289\u00f72=\", \"962\u00f76=\"],\n  [\"959\u00f73=\", \"159\u00f74=\"],\n  [\"918\u00f78=\", \"838\u00f79=\"],\n  [\"853\u00f77=\", \"234\u00f77=\"],\n  [\"763\u00f72=\", \"355\u00f72=\"],\n  [\"611\u00f74=\", \"770\u00f72=\"],\n  [\"387\u00f74=\", \"501\u00f72=\"],\n  [\"371\u00f79=\", \"609\u00f79=\"],\n  [\"701\u00f76=\", \"487\u00f77=\"],\n  [\"148\u00f73=\", \"633\u00f75=\"],\n  [\"497\u00f73=\", \"456\u00f76=\"],\n  [\"444\u00f74=\", \"684\u00f72=\"],\n  [\"359\u00f75=\", \"133\u00f76=\"],\n  [\"639\u00f74=\", \"587\u00f72=\"],\n  [\"262\u00f79=\", \"339\u00f75=\"],\n  [\"756\u00f73=\", \"791\u00f75=\"],\n  [\"768\u00f74=\", \"950\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division problems in the table with their new values.\n# Each old value is unique in the document, so Find/Replace keyed on\n# the full \"NNN\u00f7N=\" text is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"919\u00f76=\"; New = \"512\u00f75=\"},\n    @{Old = \"585\u00f79=\"; New = \"397\u00f76=\"},\n    @{Old = \"582\u00f72=\"; New = \"685\u00f74=\"},\n    @{Old = \"372\u00f74=\"; New = \"141\u00f77=\"},\n    @{Old = \"843\u00f74=\"; New = \"697\u00f79=\"},\n    @{Old = \"983\u00f74=\"; New = \"193\u00f79=\"},\n    @{Old = \"964\u00f76=\"; New = \"536\u00f72=\"},\n    @{Old = \"583\u00f73=\"; New = \"653\u00f78=\"},\n    @{Old = \"289\u00f72=\"; New = \"962\u00f76=\"},\n    @{Old = \"959\u00f73=\"; New = \"159\u00f74=\"},\n    @{Old = \"918\u00f78=\"; New = \"838\u00f79=\"},\n    @{Old = \"853\u00f77=\"; New = \"234\u00f77=\"},\n    @{Old = \"763\u00f72=\"; New = \"355\u00f72=\"},\n    @{Old = \"611\u00f74=\"; New = \"770\u00f72=\"},\n    @{Old = \"387\u00f74=\"; New = \"501\u00f72=\"},\n    @{Old = \"371\u00f79=\"; New = \"609\u00f79=\"},\n    @{Old = \"701\u00f76=\"; New = \"487\u00f77=\"},\n    @{Old = \"148\u00f73=\"; New = \"633\u00f75=\"},\n    @{Old = \"497\u00f73=\"; New = \"456\u00f76=\"},\n    @{Old = \"444\u00f74=\"; New = \"684\u00f72=\"},\n    @{Old = \"359\u00f75=\"; New = \"133\u00f76=\"},\n    @{Old = \"639\u00f74=\"; New = \"587\u00f72=\"},\n    @{Old = \"262\u00f79=\"; New = \"339\u00f75=\"},\n    @{Old = \"756\u00f73=\"; New = \"791\u00f75=\"},\n    @{Old = \"768\u00f74=\"; New = \"950\u00f77=\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2) | Out-Null\n}\n"}
